$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 679.7368  # H28: was 650.75
$ws.Cells.Item(28, 10).Value = 3006  # J28: was 1553
$ws.Cells.Item(28, 12).Value = 3006  # L28: was 1553
$ws.Cells.Item(28, 14).Value = -3976  # N28: was -2523
# Row 58
$ws.Cells.Item(58, 8).Value = 5362.5  # H58: was 5506.25
$ws.Cells.Item(58, 9).Value = 50  # I58: was 1137.5
$ws.Cells.Item(58, 10).Value = 7133.3335  # J58: was 9875
$ws.Cells.Item(58, 11).Value = 150  # K58: was 3412.5
$ws.Cells.Item(58, 12).Value = 21400.0005  # L58: was 29625
$ws.Cells.Item(58, 13).Value = 0  # M58: was -3262.5
$ws.Cells.Item(58, 14).Value = -21700.0005  # N58: was -29925
# Row 98
$ws.Cells.Item(98, 8).Value = 7204.16  # H98: was 7208.077
$ws.Cells.Item(98, 9).Value = 2919.238  # I98: was 2990.5
$ws.Cells.Item(98, 10).Value = 29700  # J98: was 21266.666
$ws.Cells.Item(98, 11).Value = 2919.238  # K98: was 2990.5
$ws.Cells.Item(98, 12).Value = 29700  # L98: was 21266.666
$ws.Cells.Item(98, 13).Value = -1421.238  # M98: was -1492.5
$ws.Cells.Item(98, 14).Value = -32696  # N98: was -24262.666
# Row 103
$ws.Cells.Item(103, 8).Value = 1509.9  # H103: was 1346.4166
$ws.Cells.Item(103, 9).Value = 2090.5  # I103: was 1867.4286
$ws.Cells.Item(103, 10).Value = 639  # J103: was 617
$ws.Cells.Item(103, 11).Value = 6271.5  # K103: was 5602.2858
$ws.Cells.Item(103, 12).Value = 1917  # L103: was 1851
$ws.Cells.Item(103, 13).Value = -5685.5  # M103: was -5016.2858
$ws.Cells.Item(103, 14).Value = -3089  # N103: was -3023
# Row 116
$ws.Cells.Item(116, 8).Value = 9611.619000000001  # H116: was 9945
$ws.Cells.Item(116, 9).Value = 16426.25  # I116: was 21235
$ws.Cells.Item(116, 10).Value = 8008.1763  # J116: was 8063.3335
$ws.Cells.Item(116, 11).Value = 16426.25  # K116: was 21235
$ws.Cells.Item(116, 12).Value = 8008.1763  # L116: was 8063.3335
$ws.Cells.Item(116, 13).Value = -12984.25  # M116: was -17793
$ws.Cells.Item(116, 14).Value = -14892.1763  # N116: was -14947.3335
# Row 122
$ws.Cells.Item(122, 8).Value = 7204.16  # H122: was 7208.077
$ws.Cells.Item(122, 9).Value = 2919.238  # I122: was 2990.5
$ws.Cells.Item(122, 10).Value = 29700  # J122: was 21266.666
$ws.Cells.Item(122, 11).Value = 8757.714  # K122: was 8971.5
$ws.Cells.Item(122, 12).Value = 89100  # L122: was 63799.99800000001
$ws.Cells.Item(122, 13).Value = -6307.714  # M122: was -6521.5
$ws.Cells.Item(122, 14).Value = -94000  # N122: was -68699.99800000001
# Row 134
$ws.Cells.Item(134, 8).Value = 43063.07  # H134: was 43048.777
$ws.Cells.Item(134, 10).Value = 43063.07  # J134: was 43048.777
$ws.Cells.Item(134, 12).Value = 43063.07  # L134: was 43048.777
$ws.Cells.Item(134, 14).Value = -53203.07  # N134: was -53188.777
# Row 137
$ws.Cells.Item(137, 8).Value = 280676.3  # H137: was 280687.38
$ws.Cells.Item(137, 10).Value = 3359.8333  # J137: was 3426.1667
$ws.Cells.Item(137, 12).Value = 10079.4999  # L137: was 10278.5001
$ws.Cells.Item(137, 14).Value = -15179.4999  # N137: was -15378.5001
# Row 138
$ws.Cells.Item(138, 8).Value = 4029.795  # H138: was 4133.946
$ws.Cells.Item(138, 9).Value = 3181.68  # I138: was 3264.25
$ws.Cells.Item(138, 10).Value = 5544.2856  # J138: was 5739.5386
$ws.Cells.Item(138, 11).Value = 9545.039999999999  # K138: was 9792.75
$ws.Cells.Item(138, 12).Value = 16632.8568  # L138: was 17218.6158
$ws.Cells.Item(138, 13).Value = -4405.039999999999  # M138: was -4652.75
$ws.Cells.Item(138, 14).Value = -26912.8568  # N138: was -27498.6158

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 46
$ws.Cells.Item(46, 8).Value = 9575.5  # H46: was 27348
$ws.Cells.Item(46, 9).Value = 4699  # I46: was 27348
$ws.Cells.Item(46, 10).Value = 14452  # J46: was 0
$ws.Cells.Item(46, 11).Value = 4699  # K46: was 27348
$ws.Cells.Item(46, 12).Value = 14452  # L46: was 0
$ws.Cells.Item(46, 13).Value = -4380  # M46: was -27029
$ws.Cells.Item(46, 14).Value = -15090  # N46: new cell
# Row 61
$ws.Cells.Item(61, 8).Value = 8519.611000000001  # H61: was 9136.294
$ws.Cells.Item(61, 9).Value = 3487.6924  # I61: was 3942
$ws.Cells.Item(61, 11).Value = 3487.6924  # K61: was 3942
$ws.Cells.Item(61, 13).Value = -3275.6924  # M61: was -3730
# Row 114
$ws.Cells.Item(114, 8).Value = 21198.5  # H114: was 24299.334
$ws.Cells.Item(114, 10).Value = 21198.5  # J114: was 24299.334
$ws.Cells.Item(114, 12).Value = 21198.5  # L114: was 24299.334
$ws.Cells.Item(114, 14).Value = -29876.5  # N114: was -32977.334
# Row 136
$ws.Cells.Item(136, 8).Value = 8519.611000000001  # H136: was 9136.294
$ws.Cells.Item(136, 9).Value = 3487.6924  # I136: was 3942
$ws.Cells.Item(136, 11).Value = 10463.0772  # K136: was 11826
$ws.Cells.Item(136, 13).Value = -7913.0772  # M136: was -9276

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Cells.Item(22, 8).Value = 345.45456  # H22: was 285.13333
$ws.Cells.Item(22, 10).Value = 807.5  # J22: was 348.66666
$ws.Cells.Item(22, 12).Value = 807.5  # L22: was 348.66666
$ws.Cells.Item(22, 14).Value = -1153.5  # N22: was -694.66666
# Row 132
$ws.Cells.Item(132, 8).Value = 0  # H132: was 127000
$ws.Cells.Item(132, 10).Value = 0  # J132: was 127000
$ws.Cells.Item(132, 12).Value = 0  # L132: was 127000
$ws.Cells.Item(132, 14).ClearContents()  # N132: was -137120
# Row 134
$ws.Cells.Item(134, 8).Value = 8365.5  # H134: was 7837.161
$ws.Cells.Item(134, 9).Value = 5735.6816  # I134: was 5395.04
$ws.Cells.Item(134, 10).Value = 18008.166  # J134: was 18012.666
$ws.Cells.Item(134, 11).Value = 17207.0448  # K134: was 16185.12
$ws.Cells.Item(134, 12).Value = 54024.49800000001  # L134: was 54037.99800000001
$ws.Cells.Item(134, 13).Value = -14672.0448  # M134: was -13650.12
$ws.Cells.Item(134, 14).Value = -59094.49800000001  # N134: was -59107.99800000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 3289.7  # H31: was 3612.375
$ws.Cells.Item(31, 9).Value = 2079.4  # I31: was 2133
$ws.Cells.Item(31, 11).Value = 2079.4  # K31: was 2133
$ws.Cells.Item(31, 13).Value = -1784.4  # M31: was -1838
# Row 34
$ws.Cells.Item(34, 8).Value = 3289.7  # H34: was 3612.375
$ws.Cells.Item(34, 9).Value = 2079.4  # I34: was 2133
$ws.Cells.Item(34, 11).Value = 2079.4  # K34: was 2133
$ws.Cells.Item(34, 13).Value = -1877.4  # M34: was -1931
# Row 58
$ws.Cells.Item(58, 8).Value = 4114.7646  # H58: was 4462.3125
$ws.Cells.Item(58, 9).Value = 2431.5833  # I58: was 2844.3333
$ws.Cells.Item(58, 10).Value = 8154.4  # J58: was 6542.5713
$ws.Cells.Item(58, 11).Value = 2431.5833  # K58: was 2844.3333
$ws.Cells.Item(58, 12).Value = 8154.4  # L58: was 6542.5713
$ws.Cells.Item(58, 13).Value = -2228.5833  # M58: was -2641.3333
$ws.Cells.Item(58, 14).Value = -8560.4  # N58: was -6948.5713
# Row 93
$ws.Cells.Item(93, 8).Value = 32679.572  # H93: was 30715.285
$ws.Cells.Item(93, 9).Value = 26459.666  # I93: was 24168
$ws.Cells.Item(93, 11).Value = 26459.666  # K93: was 24168
$ws.Cells.Item(93, 13).Value = -24587.666  # M93: was -22296
# Row 132
$ws.Cells.Item(132, 8).Value = 2561.0386  # H132: was 2864.4546
$ws.Cells.Item(132, 9).Value = 2563.48  # I132: was 2864.4546
$ws.Cells.Item(132, 10).Value = 2500  # J132: was 0
$ws.Cells.Item(132, 11).Value = 7690.440000000001  # K132: was 8593.363799999999
$ws.Cells.Item(132, 12).Value = 7500  # L132: was 0
$ws.Cells.Item(132, 13).Value = -5160.440000000001  # M132: was -6063.363799999999
$ws.Cells.Item(132, 14).Value = -12560  # N132: new cell
# Row 136
$ws.Cells.Item(136, 8).Value = 4114.7646  # H136: was 4462.3125
$ws.Cells.Item(136, 9).Value = 2431.5833  # I136: was 2844.3333
$ws.Cells.Item(136, 10).Value = 8154.4  # J136: was 6542.5713
$ws.Cells.Item(136, 11).Value = 7294.749899999999  # K136: was 8532.999899999999
$ws.Cells.Item(136, 12).Value = 24463.2  # L136: was 19627.7139
$ws.Cells.Item(136, 13).Value = -4744.749899999999  # M136: was -5982.999899999999
$ws.Cells.Item(136, 14).Value = -29563.2  # N136: was -24727.7139

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 99
$ws.Cells.Item(99, 8).Value = 5549.8335  # H99: was 6509.8
$ws.Cells.Item(99, 10).Value = 10583  # J99: was 15499.5
$ws.Cells.Item(99, 12).Value = 31749  # L99: was 46498.5
$ws.Cells.Item(99, 14).Value = -36241  # N99: was -50990.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 106
$ws.Cells.Item(106, 8).Value = 45333.668  # H106: was 46665
$ws.Cells.Item(106, 10).Value = 45333.668  # J106: was 46665
$ws.Cells.Item(106, 12).Value = 45333.668  # L106: was 46665
$ws.Cells.Item(106, 14).Value = -47857.668  # N106: was -49189

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 2300  # H7: was 0
$ws.Cells.Item(7, 9).Value = 2300  # I7: was 0
$ws.Cells.Item(7, 11).Value = 2300  # K7: was 0
$ws.Cells.Item(7, 13).Value = -2188  # M7: new cell
# Row 22
$ws.Cells.Item(22, 8).Value = 2374  # H22: was 1499.2222
$ws.Cells.Item(22, 9).Value = 1998.5  # I22: was 1142
$ws.Cells.Item(22, 11).Value = 1998.5  # K22: was 1142
$ws.Cells.Item(22, 13).Value = -1703.5  # M22: was -847
# Row 27
$ws.Cells.Item(27, 8).Value = 2374  # H27: was 1499.2222
$ws.Cells.Item(27, 9).Value = 1998.5  # I27: was 1142
$ws.Cells.Item(27, 11).Value = 1998.5  # K27: was 1142
$ws.Cells.Item(27, 13).Value = -1891.5  # M27: was -1035
# Row 76
$ws.Cells.Item(76, 8).Value = 100000  # H76: was 63999.5
$ws.Cells.Item(76, 10).Value = 100000  # J76: was 63999.5
$ws.Cells.Item(76, 12).Value = 100000  # L76: was 63999.5
$ws.Cells.Item(76, 14).Value = -100676  # N76: was -64675.5
# Row 79
$ws.Cells.Item(79, 8).Value = 100000  # H79: was 63999.5
$ws.Cells.Item(79, 10).Value = 100000  # J79: was 63999.5
$ws.Cells.Item(79, 12).Value = 100000  # L79: was 63999.5
$ws.Cells.Item(79, 14).Value = -102340  # N79: was -66339.5
# Row 126
$ws.Cells.Item(126, 8).Value = 2300  # H126: was 0
$ws.Cells.Item(126, 9).Value = 2300  # I126: was 0
$ws.Cells.Item(126, 11).Value = 6900  # K126: was 0
$ws.Cells.Item(126, 13).Value = -4430  # M126: new cell
# Row 132
$ws.Cells.Item(132, 8).Value = 3813.524  # H132: was 3917.318
$ws.Cells.Item(132, 9).Value = 3331.6667  # I132: was 3371.7334
$ws.Cells.Item(132, 10).Value = 5018.1665  # J132: was 5086.4287
$ws.Cells.Item(132, 11).Value = 9995.000100000001  # K132: was 10115.2002
$ws.Cells.Item(132, 12).Value = 15054.4995  # L132: was 15259.2861
$ws.Cells.Item(132, 13).Value = -7465.000100000001  # M132: was -7585.200199999999
$ws.Cells.Item(132, 14).Value = -20114.4995  # N132: was -20319.2861

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Cells.Item(75, 8).Value = 62559  # H75: was 46706
$ws.Cells.Item(75, 9).Value = 60118  # I75: was 46706
$ws.Cells.Item(75, 10).Value = 65000  # J75: was 0
$ws.Cells.Item(75, 11).Value = 60118  # K75: was 46706
$ws.Cells.Item(75, 12).Value = 65000  # L75: was 0
$ws.Cells.Item(75, 13).Value = -59182  # M75: was -45770
$ws.Cells.Item(75, 14).Value = -66872  # N75: new cell
# Row 78
$ws.Cells.Item(78, 8).Value = 62559  # H78: was 46706
$ws.Cells.Item(78, 9).Value = 60118  # I78: was 46706
$ws.Cells.Item(78, 10).Value = 65000  # J78: was 0
$ws.Cells.Item(78, 11).Value = 180354  # K78: was 140118
$ws.Cells.Item(78, 12).Value = 195000  # L78: was 0
$ws.Cells.Item(78, 13).Value = -175674  # M78: was -135438
$ws.Cells.Item(78, 14).Value = -204360  # N78: new cell
# Row 132
$ws.Cells.Item(132, 8).Value = 3926.3845  # H132: was 4949.875
$ws.Cells.Item(132, 9).Value = 1557.3334  # I132: was 1700
$ws.Cells.Item(132, 10).Value = 5957  # J132: was 6033.1665
$ws.Cells.Item(132, 11).Value = 4672.0002  # K132: was 5100
$ws.Cells.Item(132, 12).Value = 17871  # L132: was 18099.4995
$ws.Cells.Item(132, 13).Value = -2142.0002  # M132: was -2570
$ws.Cells.Item(132, 14).Value = -22931  # N132: was -23159.4995
# Row 136
$ws.Cells.Item(136, 8).Value = 4137.304  # H136: was 4138.0435
$ws.Cells.Item(136, 9).Value = 3746.2788  # I136: was 3708.4355
$ws.Cells.Item(136, 10).Value = 7118.875  # J136: was 7943.143
$ws.Cells.Item(136, 11).Value = 11238.8364  # K136: was 11125.3065
$ws.Cells.Item(136, 12).Value = 21356.625  # L136: was 23829.429
$ws.Cells.Item(136, 13).Value = -8688.8364  # M136: was -8575.306500000001
$ws.Cells.Item(136, 14).Value = -26456.625  # N136: was -28929.429
